# Update cryptocurrency price/volume data as scraped on Wed May 10 18:41:58 UTC 2023.
# Columns: B=Coin name, C=Link, D=Price (text), E=Volume(1h) (text, padded with spaces).
# Row 46/47 additionally swap the Cronos/Quant entries (sort-order change).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'27.520.28"
$ws.Range("E2").Value = "  -0.57%  "

# Row 3
$ws.Range("D3").Value = "'1.838.81"
$ws.Range("E3").Value = "  -0.61%  "

# Row 4
$ws.Range("E4").Value = "  +0.33%  "

# Row 5
$ws.Range("D5").Value = "'313.46"
$ws.Range("E5").Value = "  +0.12%  "

# Row 6
$ws.Range("E6").Value = "  +0.24%  "

# Row 7
$ws.Range("D7").Value = "'0.4249"
$ws.Range("E7").Value = "  +0.35%  "

# Row 8
$ws.Range("D8").Value = "'0.3661"
$ws.Range("E8").Value = "  +0.45%  "

# Row 9
$ws.Range("D9").Value = "'0.07246"
$ws.Range("E9").Value = "  -0.64%  "

# Row 10
$ws.Range("D10").Value = "'0.8683"
$ws.Range("E10").Value = "  -1.03%  "

# Row 11
$ws.Range("D11").Value = "'20.78"
$ws.Range("E11").Value = "  +0.67%  "

# Row 12
$ws.Range("D12").Value = "'1.852.20"
$ws.Range("E12").Value = "  +0.12%  "

# Row 13
$ws.Range("D13").Value = "'5.391"
$ws.Range("E13").Value = "  +1.20%  "

# Row 14
$ws.Range("D14").Value = "'6.519"
$ws.Range("E14").Value = "  -0.13%  "

# Row 15
$ws.Range("D15").Value = "'0.06940"

# Row 16
$ws.Range("D16").Value = "'1.003"
$ws.Range("E16").Value = "  +0.32%  "

# Row 17
$ws.Range("D17").Value = "'80.19"
$ws.Range("E17").Value = "  +0.60%  "

# Row 18
$ws.Range("D18").Value = "'0.000009015"
$ws.Range("E18").Value = "  +0.95%  "

# Row 19
$ws.Range("E19").Value = "  +0.19%  "

# Row 20
$ws.Range("D20").Value = "'15.45"
$ws.Range("E20").Value = "  +0.76%  "

# Row 21
$ws.Range("D21").Value = "'27.536.26"
$ws.Range("E21").Value = "  -0.55%  "

# Row 22
$ws.Range("D22").Value = "'5.059"
$ws.Range("E22").Value = "  +1.56%  "

# Row 23
$ws.Range("D23").Value = "'10.85"
$ws.Range("E23").Value = "  +4.79%  "

# Row 24
$ws.Range("D24").Value = "'2.071.11"
$ws.Range("E24").Value = "  -0.14%  "

# Row 25
$ws.Range("E25").Value = "  -1.37%  "

# Row 26
$ws.Range("D26").Value = "'154.08"
$ws.Range("E26").Value = "  -0.10%  "

# Row 27
$ws.Range("D27").Value = "'18.32"
$ws.Range("E27").Value = "  -2.73%  "

# Row 28
$ws.Range("D28").Value = "'5.241"
$ws.Range("E28").Value = "  -0.36%  "

# Row 29
$ws.Range("D29").Value = "'115.18"
$ws.Range("E29").Value = "  -5.69%  "

# Row 30
$ws.Range("D30").Value = "'1.851"
$ws.Range("E30").Value = "  -1.17%  "

# Row 31
$ws.Range("D31").Value = "'0.08870"
$ws.Range("E31").Value = "  +0.21%  "

# Row 32
$ws.Range("D32").Value = "'0.7747"
$ws.Range("E32").Value = "  +1.00%  "

# Row 33
$ws.Range("D33").Value = "'4.543"
$ws.Range("E33").Value = "  +0.01%  "

# Row 34
$ws.Range("E34").Value = "  -0.58%  "

# Row 35
$ws.Range("E35").Value = "  +4.11%  "

# Row 36
$ws.Range("E36").Value = "  +0.28%  "

# Row 37
$ws.Range("D37").Value = "'1.101"
$ws.Range("E37").Value = "  +0.55%  "

# Row 38
$ws.Range("D38").Value = "'0.05386"
$ws.Range("E38").Value = "  +0.52%  "

# Row 39
$ws.Range("E39").Value = "  +0.73%  "

# Row 40
$ws.Range("D40").Value = "'2.825"
$ws.Range("E40").Value = "  +0.21%  "

# Row 41
$ws.Range("D41").Value = "'0.5123"
$ws.Range("E41").Value = "  +0.67%  "

# Row 42
$ws.Range("D42").Value = "'0.1662"
$ws.Range("E42").Value = "  +0.88%  "

# Row 43
$ws.Range("D43").Value = "'6.737"
$ws.Range("E43").Value = "  -2.28%  "

# Row 44
$ws.Range("D44").Value = "'8.516"
$ws.Range("E44").Value = "  +2.22%  "

# Row 45
$ws.Range("E45").Value = "  +1.32%  "

# Row 46
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "'106.29"
$ws.Range("E46").Value = "  +0.71%  "

# Row 47
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.06533"
$ws.Range("E47").Value = "  -0.10%  "

# Row 48
$ws.Range("D48").Value = "'0.4703"
$ws.Range("E48").Value = "  +0.38%  "

# Row 49
$ws.Range("E49").Value = "  +0.24%  "

# Row 50
$ws.Range("D50").Value = "'1.635"
$ws.Range("E50").Value = "  +0.64%  "

# Row 51
$ws.Range("D51").Value = "'1.792"
$ws.Range("E51").Value = "  +4.12%  "
